$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value for the "Price" (D) column.
# These values must be stored as TEXT (the sheet stores numbers-as-strings),
# so we force the cell's number format to Text ("@") before assigning the
# value — otherwise Excel would auto-convert the numeric-looking string into
# a real number (and sometimes into scientific notation).
$priceUpdates = [ordered]@{
    "D2"  = "243.38"
    "D3"  = "23.06"
    "D4"  = "5.401"
    "D7"  = "6.545"
    "D8"  = "0.8115"
    "D9"  = "0.9105"
    "D10" = "0.1415"
    "D11" = "0.07438"
    "D12" = "0.03267"
    "D13" = "0.03066"
    "D14" = "0.09351"
    "D15" = "3.859"
    "D16" = "0.001561"
    "D17" = "0.04682"
    "D18" = "0.0005942"
    "D19" = "0.005929"
    "D20" = "0.004981"
    "D21" = "0.0009830"
    "D23" = "3.609"
    "D40" = "0.03960"
    "D41" = "0.006184"
    "D42" = "0.1074"
    "D44" = "0.008815"
    "D45" = "0.00005161"
    "D47" = "0.7202"
    "D48" = "0.002267"
}

foreach ($addr in $priceUpdates.Keys) {
    $range = $ws.Range($addr)
    $range.NumberFormat = "@"
    $range.Value = $priceUpdates[$addr]
}

# Volume(1h) (E) column text updates.
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
